$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 164, shifting existing rows 164-170 down to 165-171
$ws.Rows(164).Insert()

# Populate the newly inserted row 164 with the new data record
$ws.Cells.Item(164, 1).Value = 7
$ws.Cells.Item(164, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(164, 3).Value = "Ñuble"
$ws.Cells.Item(164, 4).Value = 44509
$ws.Cells.Item(164, 5).Value = 16
$ws.Cells.Item(164, 6).Value = 100112043
$ws.Cells.Item(164, 7).Value = "Pepino ensalada"
$ws.Cells.Item(164, 8).Value = "Sin especificar"
$ws.Cells.Item(164, 9).Value = "Primera"
$ws.Cells.Item(164, 10).Value = 100
$ws.Cells.Item(164, 11).Value = 9000
$ws.Cells.Item(164, 12).Value = 10000
$ws.Cells.Item(164, 13).Value = 9500
$ws.Cells.Item(164, 14).Value = "$/caja 80 unidades"
$ws.Cells.Item(164, 15).Value = "Región del Maule"
$ws.Cells.Item(164, 16).Value = 119
$ws.Cells.Item(164, 17).Value = 80
$ws.Cells.Item(164, 18).Value = "Hortaliza"
